$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")

# Row 18: formulas are mathematically unchanged (col17 / prevcol17) but become a
# shared formula group. Re-enter them so the engine groups them together.
for ($i = 1; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $prev = $cols[$i - 1]
    $ws.Range($col + "18").Formula = "=" + $col + "17/" + $prev + "17"
}

# Row 20: C20 now uses a relative reference to B20 (was absolute $B$20)
$ws.Range("C20").Formula = "=B20*C18"

# D20:O20 now chain-multiply off the previous column's row-20 result instead of
# always referencing $B$20
for ($i = 2; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $prev = $cols[$i - 1]
    $ws.Range($col + "20").Formula = "=" + $prev + "20*" + $col + "18"
}

# O20 loses its distinct right-border style and now renders like the rest of the
# row-20 block (no border)
$ws.Range("O20").Borders.Item(10).LineStyle = -4142

# Selection moves to G26
$ws.Range("G26").Select()
